$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'48.115.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.36%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.497.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.13%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'319.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.08%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'105.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.92%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -1.16%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -4.17%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'38.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.64%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'20.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.25%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0803"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.01%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -0.32%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.13%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.887.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.96%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.500.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.85%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -3.43%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'47.932.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.50%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'InternetComputer(DFINITY)"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'13.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.83%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'ImmutableX"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'2.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +8.97%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.09%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.0₃0932"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.38%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'71.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.99%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'271.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.88%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -2.16%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.05%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'25.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.66%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -4.49%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'9.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.55%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -3.56%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'34.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.87%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'49.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.01%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.08%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'19.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.51%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -2.00%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.0773"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.33%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.59%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'4.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.06%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -4.71%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'121.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.01%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'22.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.36%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -2.13%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.27%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.94%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.999.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.76%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.77%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -1.23%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'8.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.99%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -2.03%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'79.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.00%  "
$ws.Range("E51").Style = "Normal"
